$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "展览" (Exhibitions) - sheet index 1
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Row 2: a brand new event ("萤火虫动漫游戏嘉年华 x KKWORLD2024 快看漫画乐园")
# replaces what used to be the "陈张太康" KKWORLD signing event.
$ws1.Range("B2").Value = "2024-07-19"
$ws1.Range("C2").Value = "广州·萤火虫动漫游戏嘉年华 × KKWORLD2024 快看漫画乐园"
$ws1.Range("E2").Value = "2024.07.19 09:00-07.22 17:00"
$ws1.Range("F2").Value = 41383
$ws1.Range("G2").Value = "已售罄"
$ws1.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=87210"
$ws1.Range("I2").Value = "//i1.hdslb.com/bfs/openplatform/202406/DTCdOTPs1718177177472.jpeg"

# Row 3: the old "陈张太康" event (formerly row 2) shifts down into row 3,
# replacing the old "锦鲤" event content (which stays duplicated in row 4).
$ws1.Range("B3").Value = "2024-07-20"
$ws1.Range("C3").Value = "广州·KKWORLD-【陈张太康】配音演员签名内场礼包"
$ws1.Range("E3").Value = "2024.07.20 10:30-07.20 14:00"
$ws1.Range("F3").Value = 17
$ws1.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=89072"
$ws1.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202407/CTKJTCjG1720513282520.png"

# Remaining rows on this sheet only get their "想去人数" (and occasionally
# "最低票价") counters bumped.
$ws1.Range("F5").Value = 9290
$ws1.Range("F6").Value = 194
$ws1.Range("F7").Value = 777
$ws1.Range("G7").Value = 68
$ws1.Range("F8").Value = 856
$ws1.Range("F9").Value = 697
$ws1.Range("F10").Value = 196
$ws1.Range("F12").Value = 276
$ws1.Range("F13").Value = 847
$ws1.Range("F16").Value = 687
$ws1.Range("F17").Value = 293
$ws1.Range("F18").Value = 1332
$ws1.Range("F20").Value = 593
$ws1.Range("F21").Value = 673
$ws1.Range("F22").Value = 444
$ws1.Range("F23").Value = 657
$ws1.Range("F24").Value = 705
$ws1.Range("F26").Value = 40
$ws1.Range("F27").Value = 55
$ws1.Range("F28").Value = 470
$ws1.Range("F29").Value = 498
$ws1.Range("F32").Value = 910
$ws1.Range("F34").Value = 429
$ws1.Range("F37").Value = 137
$ws1.Range("F38").Value = 347
$ws1.Range("F39").Value = 1197
$ws1.Range("F40").Value = 276
$ws1.Range("F42").Value = 1197
$ws1.Range("F43").Value = 361
$ws1.Range("F45").Value = 7
$ws1.Range("F46").Value = 17
$ws1.Range("F48").Value = 38

# ------------------------------------------------------------------
# Sheet "演出" (Performances) - sheet index 2
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F11").Value = 121

# ------------------------------------------------------------------
# Sheet "本地生活" (Local life) - sheet index 3
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 1983
$ws3.Range("F3").Value = 487
$ws3.Range("F4").Value = 339

# ------------------------------------------------------------------
# Sheet "全部类型" (All types) - sheet index 4
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1983
$ws4.Range("F3").Value = 487
$ws4.Range("F10").Value = 9290
$ws4.Range("F11").Value = 194
$ws4.Range("F12").Value = 777
$ws4.Range("G12").Value = 68
$ws4.Range("F14").Value = 339
$ws4.Range("F15").Value = 856
$ws4.Range("F16").Value = 121
$ws4.Range("F17").Value = 276
$ws4.Range("F18").Value = 847
$ws4.Range("F21").Value = 293
$ws4.Range("F22").Value = 1332
$ws4.Range("F24").Value = 593
$ws4.Range("F25").Value = 444
$ws4.Range("F26").Value = 657
$ws4.Range("F27").Value = 705
$ws4.Range("F29").Value = 55
$ws4.Range("F30").Value = 470
$ws4.Range("F33").Value = 498
$ws4.Range("F36").Value = 910
$ws4.Range("F38").Value = 429
$ws4.Range("F40").Value = 137
$ws4.Range("F41").Value = 347
$ws4.Range("F42").Value = 276
$ws4.Range("F43").Value = 1197
$ws4.Range("F44").Value = 361
$ws4.Range("F47").Value = 17
